$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, pushing existing rows 196-200 down to 197-201
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new weekly record
$ws.Range("A196").Value = 4
$ws.Range("B196").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C196").Value = "Los Lagos"
$ws.Range("D196").Value = 44568
$ws.Range("E196").Value = 10
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100102
$ws.Range("H196").Value = "Cítricos"
$ws.Range("I196").Value = 100102006
$ws.Range("J196").Value = "Pomelo"
$ws.Range("K196").Value = "Start Ruby"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 200
$ws.Range("N196").Value = 11000
$ws.Range("O196").Value = 12000
$ws.Range("P196").Value = 11500
$ws.Range("Q196").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R196").Value = "Región de O'Higgins"
$ws.Range("S196").Value = 821
$ws.Range("T196").Value = 14
